{"js": "// Replace the division problems in the document body with the updated\n// operands/divisors, preserving run formatting by using Search + Replace\n// (insertText with \"Replace\") rather than rewriting whole paragraphs.\n\nconst replacements = [\n  [\"958\u00f73=\", \"694\u00f76=\"],\n  [\"686\u00f73=\", \"280\u00f78=\"],\n  [\"266\u00f79=\", \"225\u00f76=\"],\n  [\"861\u00f75=\", \"519\u00f79=\"],\n  [\"162\u00f79=\", \"710\u00f74=\"],\n  [\"828\u00f72=\", \"732\u00f76=\"],\n  [\"753\u00f79=\", \"852\u00f74=\"],\n  [\"597\u00f73=\", \"733\u00f78=\"],\n  [\"421\u00f74=\", \"575\u00f73=\"],\n  [\"431\u00f79=\", \"425\u00f78=\"],\n  [\"976\u00f76=\", \"726\u00f74=\"],\n  [\"431\u00f72=\", \"687\u00f76=\"],\n  [\"683\u00f72=\", \"952\u00f74=\"],\n  [\"893\u00f75=\", \"647\u00f74=\"],\n  [\"108\u00f75=\", \"194\u00f73=\"],\n  [\"585\u00f77=\", \"315\u00f74=\"],\n  [\"786\u00f73=\", \"620\u00f77=\"],\n  [\"661\u00f72=\", \"267\u00f73=\"],\n  [\"828\u00f73=\", \"757\u00f74=\"],\n  [\"673\u00f76=\", \"322\u00f72=\"],\n  [\"626\u00f78=\", \"884\u00f79=\"],\n  [\"246\u00f74=\", \"349\u00f77=\"],\n  [\"738\u00f73=\", \"572\u00f79=\"],\n  [\"964\u00f74=\", \"742\u00f74=\"],\n  [\"677\u00f78=\", \"488\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the document to the new operands/divisors,\n# using Word's Find/Replace so run formatting (font, size) is preserved.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"958\u00f73=\", \"694\u00f76=\"),\n    @(\"686\u00f73=\", \"280\u00f78=\"),\n    @(\"266\u00f79=\", \"225\u00f76=\"),\n    @(\"861\u00f75=\", \"519\u00f79=\"),\n    @(\"162\u00f79=\", \"710\u00f74=\"),\n    @(\"828\u00f72=\", \"732\u00f76=\"),\n    @(\"753\u00f79=\", \"852\u00f74=\"),\n    @(\"597\u00f73=\", \"733\u00f78=\"),\n    @(\"421\u00f74=\", \"575\u00f73=\"),\n    @(\"431\u00f79=\", \"425\u00f78=\"),\n    @(\"976\u00f76=\", \"726\u00f74=\"),\n    @(\"431\u00f72=\", \"687\u00f76=\"),\n    @(\"683\u00f72=\", \"952\u00f74=\"),\n    @(\"893\u00f75=\", \"647\u00f74=\"),\n    @(\"108\u00f75=\", \"194\u00f73=\"),\n    @(\"585\u00f77=\", \"315\u00f74=\"),\n    @(\"786\u00f73=\", \"620\u00f77=\"),\n    @(\"661\u00f72=\", \"267\u00f73=\"),\n    @(\"828\u00f73=\", \"757\u00f74=\"),\n    @(\"673\u00f76=\", \"322\u00f72=\"),\n    @(\"626\u00f78=\", \"884\u00f79=\"),\n    @(\"246\u00f74=\", \"349\u00f77=\"),\n    @(\"738\u00f73=\", \"572\u00f79=\"),\n    @(\"964\u00f74=\", \"742\u00f74=\"),\n    @(\"677\u00f78=\", \"488\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
